# Update cryptocurrency price/volume data in the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.571.05"
$ws.Range("E2").Value = "'  -0.07%  "
$ws.Range("D3").Value = "'1.649.06"
$ws.Range("E3").Value = "'  -0.51%  "
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("D5").Value = "'212.60"
$ws.Range("E5").Value = "'  -1.16%  "
$ws.Range("D6").Value = "'0.535"
$ws.Range("E6").Value = "'  +5.12%  "
$ws.Range("E7").Value = "'  +0.01%  "
$ws.Range("D8").Value = "'23.63"
$ws.Range("E8").Value = "'  -2.11%  "
$ws.Range("E9").Value = "'  -1.44%  "
$ws.Range("E10").Value = "'  -1.08%  "
$ws.Range("D11").Value = "'0.0891"
$ws.Range("E11").Value = "'  +1.38%  "
$ws.Range("D12").Value = "'1.884.03"
$ws.Range("E12").Value = "'  -0.53%  "
$ws.Range("D13").Value = "'1.649.45"
$ws.Range("E13").Value = "'  -0.35%  "
$ws.Range("D14").Value = "'0.586"
$ws.Range("E14").Value = "'  +3.78%  "
$ws.Range("E15").Value = "'  -2.10%  "
$ws.Range("D16").Value = "'64.59"
$ws.Range("E16").Value = "'  -1.91%  "
$ws.Range("D17").Value = "'27.556.94"
$ws.Range("E17").Value = "'  -0.03%  "
$ws.Range("D18").Value = "'231.65"
$ws.Range("E18").Value = "'  -3.68%  "
$ws.Range("E19").Value = "'  -0.34%  "
$ws.Range("D20").Value = "'7.56"
$ws.Range("E20").Value = "'  -0.66%  "
$ws.Range("E21").Value = "'  +0.01%  "
$ws.Range("D22").Value = "'4.33"
$ws.Range("E22").Value = "'  -3.51%  "
$ws.Range("D23").Value = "'9.81"
$ws.Range("E23").Value = "'  +4.76%  "
$ws.Range("D24").Value = "'2.01"
$ws.Range("E24").Value = "'  -2.27%  "
$ws.Range("D25").Value = "'148.95"
$ws.Range("E25").Value = "'  +2.13%  "
$ws.Range("E26").Value = "'  -2.60%  "
$ws.Range("E27").Value = "'  +1.81%  "
$ws.Range("E28").Value = "'  +0.09%  "
$ws.Range("D29").Value = "'15.63"
$ws.Range("E29").Value = "'  -4.09%  "
$ws.Range("E30").Value = "'  -2.51%  "
$ws.Range("E31").Value = "'  -2.89%  "
$ws.Range("E32").Value = "'  -0.52%  "
$ws.Range("D33").Value = "'3.20"
$ws.Range("E33").Value = "'  +2.77%  "
$ws.Range("D34").Value = "'1.430.34"
$ws.Range("E34").Value = "'  -2.09%  "
$ws.Range("E35").Value = "'  +1.98%  "
$ws.Range("E36").Value = "'  -0.27%  "
$ws.Range("D37").Value = "'0.569"
$ws.Range("E37").Value = "'  -0.44%  "
$ws.Range("D38").Value = "'0.887"
$ws.Range("E38").Value = "'  -4.11%  "
$ws.Range("E39").Value = "'  -3.16%  "
$ws.Range("E40").Value = "'  +0.04%  "
$ws.Range("E41").Value = "'  +0.02%  "
$ws.Range("D42").Value = "'0.818"
$ws.Range("E42").Value = "'  +3.56%  "
$ws.Range("D43").Value = "'5.54"
$ws.Range("E43").Value = "'  +2.28%  "
$ws.Range("B44").Value = "'MXToken"
$ws.Range("C44").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "'2.25"
$ws.Range("E44").Value = "'  +1.65%  "
$ws.Range("B45").Value = "'Aave"
$ws.Range("C45").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'65.29"
$ws.Range("E45").Value = "'  -6.68%  "
$ws.Range("B46").Value = "'RocketPoolETH"
$ws.Range("C46").Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "'1.792.75"
$ws.Range("E46").Value = "'  -0.52%  "
$ws.Range("B47").Value = "'RenderToken"
$ws.Range("C47").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'1.69"
$ws.Range("E47").Value = "'  -1.02%  "
$ws.Range("B48").Value = "'Quant"
$ws.Range("C48").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'88.05"
$ws.Range("E48").Value = "'  -0.60%  "
$ws.Range("B49").Value = "'BabyDogeCoin"
$ws.Range("C49").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.0₆0108"
$ws.Range("E49").Value = "'  +0.09%  "
$ws.Range("B50").Value = "'Algorand"
$ws.Range("C50").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.0993"
$ws.Range("E50").Value = "'  -3.26%  "
$ws.Range("D51").Value = "'7.77"
$ws.Range("E51").Value = "'  -0.68%  "
